$wb = $excel.ActiveWorkbook

# --- Sheet 1: "amplifiers" ---
$ws1 = $wb.Worksheets.Item(1)

# Sort the data rows (A2:G10) ascending by column A (numbers sort before text,
# matching Excel's default ascending sort behaviour).
$rng1 = $ws1.Range("A2:G10")
$rng1.Sort($ws1.Range("A2:A10"), 1)

# The row that used to hold the text label "9001" (quoted) now stores the
# literal number 9001.
$ws1.Range("A3").Value = 9001

# Sheet1 becomes the active tab, with A8 selected.
$ws1.Range("A8").Select()

# --- Sheet 2: "speakers" ---
$ws2 = $wb.Worksheets.Item(2)

# Sort the data rows (A2:F17) ascending by column A.
$rng2 = $ws2.Range("A2:F17")
$rng2.Sort($ws2.Range("A2:A17"), 1)

# speakers keeps C20 as its last remembered selection (not the active sheet).
$ws2.Range("C20").Select()

# Sheet1 ("amplifiers") is the active/selected sheet.
$ws1.Activate()
$ws1.Select()
